$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-26 Sunday" "2024-05-27 Monday"

Replace-Text "826÷3=275, 1" "344÷4=86, 0"
Replace-Text "580÷6=96, 4" "153÷3=51, 0"
Replace-Text "246÷4=61, 2" "858÷7=122, 4"
Replace-Text "748÷8=93, 4" "222÷8=27, 6"
Replace-Text "716÷5=143, 1" "759÷2=379, 1"

Replace-Text "767÷8=95, 7" "162÷6=27, 0"
Replace-Text "632÷4=158, 0" "225÷9=25, 0"
Replace-Text "370÷9=41, 1" "534÷8=66, 6"
Replace-Text "917÷8=114, 5" "675÷7=96, 3"
Replace-Text "590÷9=65, 5" "173÷2=86, 1"

Replace-Text "838÷8=104, 6" "525÷7=75, 0"
Replace-Text "167÷4=41, 3" "525÷8=65, 5"
Replace-Text "982÷6=163, 4" "910÷3=303, 1"
Replace-Text "707÷6=117, 5" "998÷7=142, 4"
Replace-Text "439÷9=48, 7" "562÷2=281, 0"

Replace-Text "217÷4=54, 1" "568÷2=284, 0"
Replace-Text "221÷8=27, 5" "343÷6=57, 1"
Replace-Text "447÷7=63, 6" "251÷6=41, 5"
Replace-Text "157÷7=22, 3" "298÷5=59, 3"
Replace-Text "607÷7=86, 5" "736÷9=81, 7"

Replace-Text "323÷8=40, 3" "945÷7=135, 0"
Replace-Text "498÷3=166, 0" "495÷8=61, 7"
Replace-Text "607÷5=121, 2" "748÷8=93, 4"
Replace-Text "443÷9=49, 2" "702÷9=78, 0"
Replace-Text "321÷2=160, 1" "784÷6=130, 4"
